$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at W, shifting the old "Pohjakoulutus maa (toinen aste)"
# column (and everything to its right) one place over to X.
$ws.Columns("W").Insert()

# Populate the new column: header + the single data row's value.
$ws.Range("W1").Value = "Toisen asteen pohjakoulutus suoritettu"
$ws.Range("W2").Value = "Kyllä"
